$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1470.7778
$ws.Range("I28").Value = 1470.7778
$ws.Range("K28").Value = 1470.7778
$ws.Range("M28").Value = -985.7778000000001

$ws.Range("H40").Value = 2000
$ws.Range("J40").Value = 2000
$ws.Range("L40").Value = 2000
$ws.Range("N40").Value = -2350

$ws.Range("H132").Value = 3046.4285
$ws.Range("I132").Value = 2042.1818
$ws.Range("K132").Value = 6126.5454
$ws.Range("M132").Value = -3596.5454

$ws.Range("H137").Value = 1502.909
$ws.Range("I137").Value = 1502.909
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 4508.727000000001
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -1958.727000000001
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3052.08
$ws.Range("I32").Value = 3024
$ws.Range("K32").Value = 3024
$ws.Range("M32").Value = -2737

$ws.Range("H61").Value = 5442.222
$ws.Range("I61").Value = 5747.5
$ws.Range("K61").Value = 5747.5
$ws.Range("M61").Value = -5535.5

$ws.Range("H110").Value = 743.44446
$ws.Range("I110").Value = 748.875
$ws.Range("K110").Value = 748.875
$ws.Range("M110").Value = 1296.125

$ws.Range("H122").Value = 5225.2173
$ws.Range("I122").Value = 5475.1
$ws.Range("K122").Value = 16425.3
$ws.Range("M122").Value = -13975.3

$ws.Range("H132").Value = 3326
$ws.Range("I132").Value = 2622
$ws.Range("J132").Value = 4499.3335
$ws.Range("K132").Value = 7866
$ws.Range("L132").Value = 13498.0005
$ws.Range("M132").Value = -5336
$ws.Range("N132").Value = -18558.0005

$ws.Range("H136").Value = 5442.222
$ws.Range("I136").Value = 5747.5
$ws.Range("K136").Value = 17242.5
$ws.Range("M136").Value = -14692.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4911.4
$ws.Range("I134").Value = 4957.1113
$ws.Range("K134").Value = 14871.3339
$ws.Range("M134").Value = -12336.3339

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 7327.6
$ws.Range("I16").Value = 7409.625
$ws.Range("J16").Value = 6999.5
$ws.Range("K16").Value = 7409.625
$ws.Range("L16").Value = 6999.5
$ws.Range("M16").Value = -7122.625
$ws.Range("N16").Value = -7573.5

$ws.Range("H31").Value = 3004
$ws.Range("I31").Value = 3004
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 3004
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -2709
$ws.Range("N31").ClearContents()

$ws.Range("H34").Value = 3004
$ws.Range("I34").Value = 3004
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 3004
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -2802
$ws.Range("N34").ClearContents()

$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

$ws.Range("H99").Value = 3853.3333
$ws.Range("I99").Value = 4021.2
$ws.Range("K99").Value = 4021.2
$ws.Range("M99").Value = -2523.2

$ws.Range("H105").Value = 3761.7144
$ws.Range("I105").Value = 2839.75
$ws.Range("J105").Value = 4991
$ws.Range("K105").Value = 2839.75
$ws.Range("L105").Value = 4991
$ws.Range("M105").Value = -1092.75
$ws.Range("N105").Value = -8485

$ws.Range("H107").Value = 2712.25
$ws.Range("I107").Value = 2000
$ws.Range("K107").Value = 2000
$ws.Range("M107").Value = -80

$ws.Range("H113").Value = 7327.6
$ws.Range("I113").Value = 7409.625
$ws.Range("J113").Value = 6999.5
$ws.Range("K113").Value = 7409.625
$ws.Range("L113").Value = 6999.5
$ws.Range("M113").Value = -5239.625
$ws.Range("N113").Value = -11339.5

$ws.Range("H122").Value = 2222
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 2222
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 6666
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -11566

$ws.Range("H126").Value = 3853.3333
$ws.Range("I126").Value = 4021.2
$ws.Range("K126").Value = 12063.6
$ws.Range("M126").Value = -9593.599999999999

$ws.Range("H134").Value = 4954.2144
$ws.Range("I134").Value = 4905.25
$ws.Range("K134").Value = 14715.75
$ws.Range("M134").Value = -12180.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

$ws.Range("H122").Value = 898
$ws.Range("I122").Value = 579.3333
$ws.Range("J122").Value = 1034.5714
$ws.Range("K122").Value = 5213.9997
$ws.Range("L122").Value = 9311.142600000001
$ws.Range("M122").Value = -2763.9997
$ws.Range("N122").Value = -14211.1426

$ws.Range("H131").Value = 1607
$ws.Range("J131").Value = 1750
$ws.Range("L131").Value = 5250
$ws.Range("N131").Value = -15330

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2475.5334
$ws.Range("I97").Value = 594.6667
$ws.Range("K97").Value = 594.6667
$ws.Range("M97").Value = -98.66669999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1534.9286
$ws.Range("I16").Value = 1419.0834
$ws.Range("K16").Value = 1419.0834
$ws.Range("M16").Value = -1249.0834

$ws.Range("H22").Value = 1463.3334
$ws.Range("I22").Value = 1490
$ws.Range("J22").Value = 1450
$ws.Range("K22").Value = 1490
$ws.Range("L22").Value = 1450
$ws.Range("M22").Value = -1195
$ws.Range("N22").Value = -2040

$ws.Range("H27").Value = 1463.3334
$ws.Range("I27").Value = 1490
$ws.Range("J27").Value = 1450
$ws.Range("K27").Value = 1490
$ws.Range("L27").Value = 1450
$ws.Range("M27").Value = -1383
$ws.Range("N27").Value = -1664

$ws.Range("H46").Value = 5000
$ws.Range("J46").Value = 5000
$ws.Range("L46").Value = 5000
$ws.Range("N46").Value = -5376

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3173.25
$ws.Range("I122").Value = 2899.3333
$ws.Range("K122").Value = 8697.999899999999
$ws.Range("M122").Value = -6247.999899999999

$ws.Range("H136").Value = 9812.08
$ws.Range("I136").Value = 9823
$ws.Range("J136").Value = 9686.5
$ws.Range("K136").Value = 29469
$ws.Range("L136").Value = 29059.5
$ws.Range("M136").Value = -26919
$ws.Range("N136").Value = -34159.5
